$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Danish"

$ws.Range("A2").Value = "UK will ban mobile carriers from selling locked handsets in 2021"
$ws.Range("B2").Value = "Danish"
$ws.Range("C2").Value = "Storbritannien forbyder mobiloperatører at sælge låste håndsæt i 2021"
$ws.Range("D2").Value = "https://www.engadget.com/uk-bans-mobile-companies-from-selling-locked-handsets-103533257.html"

$ws.Range("A3").Value = "Tesco fast tracks emissions target with solar push - Reuters UK"
$ws.Range("B3").Value = "Danish"
$ws.Range("C3").Value = "Tesco sporer hurtigt emissionsmål med solskub - Reuters UK"
$ws.Range("D3").Value = "https://uk.reuters.com/article/uk-tesco-emissions-idUKKBN27T001"

$ws.Range("A4").Value = "Asda, Tesco, Sainsbury's and Aldi warn shoppers what they can buy during second lockdown - Kent Live"
$ws.Range("B4").Value = "Danish"
$ws.Range("C4").Value = "Asda, Tesco, Sainsbury og Aldi advarer shoppere om, hvad de kan købe under anden lockdown - Kent Live"
$ws.Range("D4").Value = "https://www.kentlive.news/whats-on/shopping/asda-tesco-sainsburys-aldi-warn-4680621"
